$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 2)
$c.Value = "'Bitcoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 3)
$c.Value = "'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 4)
$c.Value = "'42.888.47"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +4.18%  "

$c = $ws.Cells.Item(3, 2)
$c.Value = "'Ethereum"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 3)
$c.Value = "'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.Value = "'2.283.85"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +5.06%  "

$c = $ws.Cells.Item(4, 2)
$c.Value = "'TetherUSD"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 3)
$c.Value = "'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.07%  "

$c = $ws.Cells.Item(5, 2)
$c.Value = "'BNB"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 3)
$c.Value = "'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.Value = "'252.52"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.99%  "

$c = $ws.Cells.Item(6, 2)
$c.Value = "'XRP"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 3)
$c.Value = "'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.Value = "'0.640"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +4.51%  "

$c = $ws.Cells.Item(7, 2)
$c.Value = "'Solana"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 3)
$c.Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 4)
$c.Value = "'72.79"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +10.12%  "

$c = $ws.Cells.Item(8, 2)
$c.Value = "'USDC"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 3)
$c.Value = "'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.12%  "

$c = $ws.Cells.Item(9, 2)
$c.Value = "'Cardano"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 3)
$c.Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.657"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +13.91%  "

$c = $ws.Cells.Item(10, 2)
$c.Value = "'Avalanche"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 3)
$c.Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 4)
$c.Value = "'38.81"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +7.22%  "

$c = $ws.Cells.Item(11, 2)
$c.Value = "'Dogecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 3)
$c.Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.0979"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +4.94%  "

$c = $ws.Cells.Item(12, 2)
$c.Value = "'OKB"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 3)
$c.Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.Value = "'59.74"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.55%  "

$c = $ws.Cells.Item(13, 2)
$c.Value = "'Polkadot"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 3)
$c.Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.Value = "'7.40"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +8.50%  "

$c = $ws.Cells.Item(14, 2)
$c.Value = "'TRON"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 3)
$c.Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.Value = "'0.106"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.14%  "

$c = $ws.Cells.Item(15, 2)
$c.Value = "'WrappedliquidstakedEther2.0"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 3)
$c.Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.Value = "'2.620.81"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +4.73%  "

$c = $ws.Cells.Item(16, 2)
$c.Value = "'Chainlink"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 3)
$c.Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.Value = "'15.04"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +5.54%  "

$c = $ws.Cells.Item(17, 2)
$c.Value = "'Polygon"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 3)
$c.Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.Value = "'0.891"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +5.34%  "

$c = $ws.Cells.Item(18, 2)
$c.Value = "'WrappedEther"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 3)
$c.Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.Value = "'2.289.81"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +6.25%  "

$c = $ws.Cells.Item(19, 2)
$c.Value = "'WrappedBTC"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 3)
$c.Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.Value = "'42.797.67"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +4.12%  "

$c = $ws.Cells.Item(20, 2)
$c.Value = "'ShibaInu"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 3)
$c.Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.Value = "'0.0000101"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +7.25%  "

$c = $ws.Cells.Item(21, 2)
$c.Value = "'Uniswap"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 3)
$c.Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.Value = "'6.36"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +5.32%  "

$c = $ws.Cells.Item(22, 2)
$c.Value = "'Litecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 3)
$c.Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.Value = "'73.53"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.71%  "

$c = $ws.Cells.Item(23, 2)
$c.Value = "'BitcoinCash"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 3)
$c.Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.Value = "'237.31"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +3.29%  "

$c = $ws.Cells.Item(24, 4)
$c.Value = "'2.15"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +6.56%  "

$c = $ws.Cells.Item(25, 4)
$c.Value = "'3.87"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.93%  "

$c = $ws.Cells.Item(26, 4)
$c.Value = "'11.69"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +2.26%  "

$c = $ws.Cells.Item(27, 4)
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.15%  "

$c = $ws.Cells.Item(28, 4)
$c.Value = "'2.45"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.88%  "

$c = $ws.Cells.Item(29, 4)
$c.Value = "'3.68"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.94%  "

$c = $ws.Cells.Item(30, 4)
$c.Value = "'2.13"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +5.41%  "

$c = $ws.Cells.Item(31, 4)
$c.Value = "'167.88"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.29%  "

$c = $ws.Cells.Item(32, 4)
$c.Value = "'21.15"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +4.95%  "

$c = $ws.Cells.Item(33, 4)
$c.Value = "'6.36"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +10.40%  "

$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.129"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +5.76%  "

$c = $ws.Cells.Item(35, 4)
$c.Value = "'0.0816"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +9.02%  "

$c = $ws.Cells.Item(36, 4)
$c.Value = "'31.07"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +27.84%  "

$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.127"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +4.89%  "

$c = $ws.Cells.Item(38, 4)
$c.Value = "'4.72"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +20.25%  "

$c = $ws.Cells.Item(39, 4)
$c.Value = "'4.79"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +6.62%  "

$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.0310"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.88%  "

$c = $ws.Cells.Item(41, 2)
$c.Value = "'Celestia"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 3)
$c.Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.Value = "'13.39"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +18.40%  "

$c = $ws.Cells.Item(42, 2)
$c.Value = "'LidoDAOToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 3)
$c.Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.Value = "'2.34"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +5.83%  "

$c = $ws.Cells.Item(43, 4)
$c.Value = "'5.99"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +9.40%  "

$c = $ws.Cells.Item(44, 4)
$c.Value = "'0.213"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +12.92%  "

$c = $ws.Cells.Item(45, 4)
$c.Value = "'9.22"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +8.67%  "

$c = $ws.Cells.Item(46, 4)
$c.Value = "'5.01"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -6.75%  "

$c = $ws.Cells.Item(47, 4)
$c.Value = "'61.45"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.57%  "

$c = $ws.Cells.Item(48, 4)
$c.Value = "'0.104"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +4.47%  "

$c = $ws.Cells.Item(49, 4)
$c.Value = "'1.19"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +4.55%  "

$c = $ws.Cells.Item(50, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.10%  "

$c = $ws.Cells.Item(51, 4)
$c.Value = "'1.20"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +5.32%  "
